$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Text content (shared strings) ----
$tabname    = "TabName"
$query_h    = "query"
$statquery_h = "StatQuery"
$dbexcel_h  = "dbExcel"
$webexcel_h = "WebExcel"

$cases_tab   = "CasesTab"
$samples_tab = "SamplesTab"
$files_tab   = "FilesTab"

$neo4j_file = "TC03_Canine_Filter_Breed-AusShephd_Neo4jData.xlsx"
$web_file   = "TC03_Canine_Filter_Breed-AusShephd_WebData.xlsx"

$stat_query = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Australian Shepherd']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$samples_query = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN ['Australian Shepherd']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, 
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS ``Sample Site``,
        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,
        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,
        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,
        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,
        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,
        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,
        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``"

$cases_query = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Australian Shepherd']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``,
        coalesce(co.cohort_description, '') AS ``Cohort``
"

$files_query = "MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Australian Shepherd']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
        coalesce(f.file_type, '') AS ``File Type``,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

# ---- Clear any pre-existing wrap formatting before re-laying out ----
$ws.Range("A1:E4").Style = "Normal"

# ---- Row 1: headers ----
$ws.Range("A1").Value = $tabname
$ws.Range("B1").Value = $query_h
$ws.Range("C1").Value = $statquery_h
$ws.Range("D1").Value = $dbexcel_h
$ws.Range("E1").Value = $webexcel_h

# ---- Row 2: Cases tab ----
$ws.Range("A2").Value = $cases_tab
$ws.Range("B2").Value = $cases_query
$ws.Range("C2").Value = $stat_query
$ws.Range("D2").Value = $neo4j_file
$ws.Range("E2").Value = $web_file

# ---- Row 3: Samples tab ----
$ws.Range("A3").Value = $samples_tab
$ws.Range("B3").Value = $samples_query
$ws.Range("C3").Value = $stat_query
$ws.Range("D3").Value = $neo4j_file
$ws.Range("E3").Value = $web_file

# ---- Row 4: Files tab ----
$ws.Range("A4").Value = $files_tab
$ws.Range("B4").Value = $files_query
$ws.Range("C4").Value = $stat_query
$ws.Range("D4").Value = $neo4j_file
$ws.Range("E4").Value = $web_file

# ---- Wrap text formatting on the query columns (B & C) for rows 2-4 ----
$ws.Range("B2:C4").WrapText = $true

# ---- Row heights (Excel auto-fit equivalents for the wrapped text) ----
$ws.Rows.Item(2).RowHeight = 275.5
$ws.Rows.Item(3).RowHeight = 232
$ws.Rows.Item(4).RowHeight = 246.5

# ---- Column widths (closest achievable values given the engine's pixel
#      quantization of ColumnWidth; targets are 10.90625 / 75.81640625 /
#      89.1796875 / 70.26953125 / 49.81640625 characters) ----
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 75
$ws.Columns.Item(3).ColumnWidth = 88.3
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 49

# ---- Sheet view: zoom + top-left cell + selection ----
$win = $excel.ActiveWindow
$win.Zoom = 55
$win.ScrollRow = 4
$win.ScrollColumn = 1
[void]$ws.Range("D12").Select()
